# Session 28 report update: regenerate MDAyuda ticket statistics report
# with fresh "Generado" timestamp and updated figures; remove the
# "Filtro aplicado" line and extend several breakdown tables with more rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Resumen General" -----------------------------------------
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Resumen General")

# Update the generation timestamp.
$ws1.Range("A2").Value = "Generado: 7 de enero de 2026, 18:03"

# Remove the "Filtro aplicado: Sistema de Ventas" row entirely; this
# shifts every row below it up by one.
$ws1.Rows(3).Delete()

# Refresh the ticket-status counts (rows shifted up after the delete).
$ws1.Range("B6").Value = 156   # Total de Tickets
$ws1.Range("B7").Value = 42    # Abiertos
$ws1.Range("B8").Value = 35    # En Proceso
$ws1.Range("B9").Value = 18    # En Espera
$ws1.Range("B10").Value = 48   # Resueltos
$ws1.Range("B11").Value = 13   # Cerrados

# ---------------------------------------------------------------------
# Sheet 2: "Por Categoria" --------------------------------------------
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Por Categoria")

$ws2.Range("A4").Value = "Portal Web"
$ws2.Range("B4").Value = 32
$ws2.Range("C4").NumberFormat = "@"
$ws2.Range("C4").Value = "20.5%"

$ws2.Range("A5").Value = "Aplicacion Movil"
$ws2.Range("B5").Value = 28
$ws2.Range("C5").NumberFormat = "@"
$ws2.Range("C5").Value = "17.9%"

$ws2.Range("A6").Value = "Sistema de Inventario"
$ws2.Range("B6").Value = 24
$ws2.Range("C6").NumberFormat = "@"
$ws2.Range("C6").Value = "15.4%"

$ws2.Range("A7").Value = "Facturacion"
$ws2.Range("B7").Value = 15
$ws2.Range("C7").NumberFormat = "@"
$ws2.Range("C7").Value = "9.6%"

$ws2.Range("A8").Value = "Otros"
$ws2.Range("B8").Value = 12
$ws2.Range("C8").NumberFormat = "@"
$ws2.Range("C8").Value = "7.7%"

# ---------------------------------------------------------------------
# Sheet 3: "Por Prioridad" --------------------------------------------
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Por Prioridad")

$ws3.Range("B3").Value = 38  # Alta
$ws3.Range("B4").Value = 78  # Media
$ws3.Range("B5").Value = 40  # Baja

# ---------------------------------------------------------------------
# Sheet 4: "Rendimiento Empleados" ------------------------------------
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rendimiento Empleados")

$ws4.Range("A4").Value = "Carlos Rodriguez"
$ws4.Range("B4").Value = 38
$ws4.Range("C4").Value = "5.1 horas"
$ws4.Range("D4").NumberFormat = "@"
$ws4.Range("D4").Value = "92%"

$ws4.Range("A5").Value = "Ana Martinez"
$ws4.Range("B5").Value = 32
$ws4.Range("C5").Value = "3.8 horas"
$ws4.Range("D5").NumberFormat = "@"
$ws4.Range("D5").Value = "98%"

$ws4.Range("A6").Value = "Luis Fernandez"
$ws4.Range("B6").Value = 28
$ws4.Range("C6").Value = "6.2 horas"
$ws4.Range("D6").NumberFormat = "@"
$ws4.Range("D6").Value = "88%"

$ws4.Range("A7").Value = "Sofia Torres"
$ws4.Range("B7").Value = 25
$ws4.Range("C7").Value = "4.5 horas"
$ws4.Range("D7").NumberFormat = "@"
$ws4.Range("D7").Value = "91%"

# ---------------------------------------------------------------------
# Sheet 5: "Cumplimiento SLA" ------------------------------------------
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Cumplimiento SLA")

$ws5.Range("B3").Value = 142  # Tickets con SLA Cumplido
$ws5.Range("B4").Value = 14   # Tickets con SLA Incumplido

Write-Host "Report regenerated."
